# Applies the "MIGO transaction processing" data refresh to Planilha1:
#   - P2/Q2 and P3/Q3 get updated purchase-order / material-document numbers
#   - R2 picks up a new (underlined-font + centered) style, matching the
#     workbook's freshly-added 8th cellXfs entry
#   - the active selection moves from Q3 to R2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: MIGO / purchase order numbers -------------------------------
$ws.Range("P2").Value = 4600244328
$ws.Range("Q2").Value = 4503342111

# --- Row 3: MIGO / purchase order numbers -------------------------------
$ws.Range("P3").Value = 4600244329
$ws.Range("Q3").Value = 4503342112

# --- R2 gets the new style (font with underline, centered horizontally) -
# This reuses fontId 2 (already underlined in the workbook) and adds a new
# cellXfs entry with centered alignment, exactly like the target diff.
$ws.Range("R2").Font.Underline = $true
$ws.Range("R2").HorizontalAlignment = -4108  # xlCenter

# --- Move the active window/selection to R2 (previously Q3) -------------
$win = $excel.ActiveWindow
$win.ScrollColumn = 8
$win.ScrollRow = 1
$ws.Range("R2").Select() | Out-Null
